$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Metadata" (sheet1): update Version, Date and Contact values
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "2.0.0"
$wsMeta.Range("B8").Value = "2024-06-03T10:45:43+02:00"
$wsMeta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# ---------------------------------------------------------------------------
# Sheet "Include from FSIII" (sheet2): add 12 new concept rows holding UUIDs,
# push the existing I-codes down below them, and append the System URI row
# ---------------------------------------------------------------------------
$wsInc = $wb.Worksheets.Item("Include from FSIII")

# Capture the existing I-code labels (currently in A2:A13) before they get
# shifted down, so we can re-write them further down the sheet afterwards.
$iCodes = @()
for ($r = 2; $r -le 13; $r++) {
    $iCodes += $wsInc.Cells.Item($r, 1).Value2
}

# Insert 12 fresh rows right before the old row 14 (the blank separator row).
# This shifts the old row 14 ("" / "") down to row 26 and the old row 15
# ("System URI" / oid) down to row 27 *without* touching their cell content,
# so those two rows keep their original representation exactly.
$wsInc.Rows("14:25").Insert()

$newUuids = @(
    "94e9c867-fbc8-4d35-8596-e6b8765b12e8",
    "55670b1e-7a36-46b2-8712-b7536237f22d",
    "9162d29a-1c7f-4585-8145-8fb4f1a999e3",
    "fa6aa904-d06e-4029-b4c4-13ead04ace27",
    "3f00a76f-8e7b-4b13-80cc-f2ceef4e51d1",
    "01150cdb-6098-48ce-bb61-60967f6bcc37",
    "1bb534f3-e526-41a9-b9c3-6157ea19c915",
    "cc377732-7f14-49b7-8940-1aa07b8884e7",
    "25dcedb3-7149-4ef9-a2c3-be30267441fb",
    "045fa500-35b0-46b7-97dd-adb60888a8ea",
    "8c539fd9-7f31-4b4e-8b30-8298c8ab640f",
    "5bfe4bda-2358-41da-946e-1fdaa33d5fe8"
)

# Rows 2-13: replace the I-code labels with the new UUID concept codes; the
# Value column (B) stays empty.
for ($i = 0; $i -lt $newUuids.Length; $i++) {
    $row = 2 + $i
    $wsInc.Cells.Item($row, 1).Value = $newUuids[$i]
    $wsInc.Cells.Item($row, 2).ClearContents()
}

# Rows 14-25 (newly inserted): the I-code labels that used to live in rows
# 2-13, with an empty Value column.
for ($i = 0; $i -lt $iCodes.Length; $i++) {
    $row = 14 + $i
    $wsInc.Cells.Item($row, 1).Value = $iCodes[$i]
    $wsInc.Cells.Item($row, 2).ClearContents()
}
